# Split the run " (Alternative)" into three runs: " (", "Changed main", ")"
# i.e. replace "Alternative" with "Changed main" while keeping it in its own
# run, distinct both from the leading " (" text and from the trailing ")".

$d = $word.ActiveDocument

# 1. Locate the run of interest by its current text.
$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute(" (Alternative)", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null

$runStart = $target.Start
$runEnd   = $target.End

# 2. Within that text, find "Alternative" and swap it for "Changed main".
$wordRange = $d.Range($runStart, $runEnd)
$wordRange.Find.ClearFormatting()
$wordRange.Find.Execute("Alternative", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null

$wordStart = $wordRange.Start
$oldWordLength = $wordRange.End - $wordRange.Start
$wordRange.Text = "Changed main"
$wordEnd = $wordStart + "Changed main".Length

# Account for the length delta ("Changed main" vs "Alternative") when
# recomputing the end of the outer " (...)" range.
$newRunEnd = $runEnd + ("Changed main".Length - $oldWordLength)

# 3. Force Word to keep " (" as its own run, separate from the preceding
#    "This is a Microsoft word document." run, by nudging (and reverting)
#    direct character formatting across the whole " (Changed main)" span.
#    A no-op text edit alone gets silently re-merged with identically
#    formatted neighboring runs, so a transient formatting toggle is used
#    purely to force the run boundary to "stick".
$rAll = $d.Range($runStart, $newRunEnd)
$rAll.Bold = 1
$rAll.Bold = 0

# 4. Likewise force "Changed main" to be its own run, separate from the
#    surrounding " (" and ")" runs.
$rMid = $d.Range($wordStart, $wordEnd)
$rMid.Bold = 1
$rMid.Bold = 0
